# Generate Report for Handoff
# The "b.md" source file has been re-handed-off for localization, so its
# status flips from "Handed back: in sync with en-US" to "Ready for
# handoff" on the Overview sheet, and the per-language (zh-cn / de-de)
# detail sheets get a fresh handoff file/timestamp plus an error detail
# explaining that the previous handback is stale.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet - row 3 is the "b.md" row
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-25 02:37:26"

# ---------------------------------------------------------------------
# zh-cn sheet - row 3 is the "b.md" row
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-25 02:37:21"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f14c03a1bd527d34d47cb6b499219a7cef99aed/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ac0adbc4ad76377a970e7dd899e169901e96430/e2e/b.md."
# Column P (Error Detail) needs to be wide enough to show the message
$zhcn.Columns.Item(16).ColumnWidth = 39.16666666666666

# ---------------------------------------------------------------------
# de-de sheet - row 3 is the "b.md" row
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-25 02:37:26"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f14c03a1bd527d34d47cb6b499219a7cef99aed/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ac0adbc4ad76377a970e7dd899e169901e96430/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.16666666666666
